# Update scripts with new TPM values (NATMI LR-pairs output regenerated).
# The "Sending cluster" assignment for the FAPs/MuSCs groups shifted (FAPs -> ECs,
# MuSCs -> FAPs) and all the expression/specificity metrics were recomputed with
# the new TPM numbers. "Target cluster" (column D) keeps the same cluster names
# per row; only columns A and E:T change value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data in column order: A, B, C, D, E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
$rows = @{
    2 = @("ECs","Fgf22","Fgfr2","ECs",1,0.3333333333333333,0.028883,0.086649,0.2501761214025038,0.2501761214025038,2,0.6666666666666666,0.05601,0.16803,0.02710547761971223,0.02710547761971223,0.00161773683,0.01455963147,0.006781143259661978,0.006781143259661978)
    3 = @("ECs","Fgf22","Fgfr2","FAPs",1,0.3333333333333333,0.028883,0.086649,0.2501761214025038,0.2501761214025038,3,1,1.864751,5.594253,0.902427539668559,0.9024275396685592,0.05385960313300001,0.484736428197,0.2257658217210843,0.2257658217210843)
    4 = @("ECs","Fgf22","Fgfr2","MuSCs",1,0.3333333333333333,0.028883,0.086649,0.2501761214025038,0.2501761214025038,3,1,0.145611,0.436833,0.07046698271172858,0.07046698271172858,0.004205682513000001,0.037851142617,0.01762915642175755,0.01762915642175755)
    5 = @("FAPs","Fgf22","Fgfr2","ECs",1,0.3333333333333333,0.08656766666666667,0.259703,0.7498238785974961,0.7498238785974962,2,0.6666666666666666,0.05601,0.16803,0.02710547761971223,0.02710547761971223,0.00484865501,0.04363789509,0.02032433436005025,0.02032433436005026)
    6 = @("FAPs","Fgf22","Fgfr2","FAPs",1,0.3333333333333333,0.08656766666666667,0.259703,0.7498238785974961,0.7498238785974962,3,1,1.864751,5.594253,0.902427539668559,0.9024275396685592,0.1614271429843333,1.452844286859,0.6766617179474748,0.6766617179474749)
    7 = @("FAPs","Fgf22","Fgfr2","MuSCs",1,0.3333333333333333,0.08656766666666667,0.259703,0.7498238785974961,0.7498238785974962,3,1,0.145611,0.436833,0.07046698271172858,0.07046698271172858,0.012605204511,0.113446840599,0.05283782628997102,0.05283782628997103)
}

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    $arr = New-Object "object[,]" 1,20
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $ws.Range("A$r`:T$r").Value = $arr
}
